$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace row 2 data (old Catrice Kimball record -> new James Calloway record)
$ws.Range("A2").Value = "James"
$ws.Range("B2").Value = "Calloway"
$ws.Range("C2").Value = "jamescalloway402@gmail.com"
$ws.Range("D2").Value = "whiqufiogheqkdvw"
$ws.Range("E2").Value = "185.14.97.29:4001"
$ws.Range("F2").Value = "eVLjgYYsF64zW8Zx"
$ws.Range("G2").Value = "gV9XKZj2J47VwkeY"

# Replace row 3 data (old Xuan Warren record -> new David Gibbons record)
$ws.Range("A3").Value = "David"
$ws.Range("B3").Value = "Gibbons"
$ws.Range("C3").Value = "0101dave.gibbons22@gmail.com"
$ws.Range("D3").Value = "wlmgbnezgtgglefg"
$ws.Range("E3").Value = "81.28.96.131:58065"
$ws.Range("F3").Value = "Ks4BcUQh43z8AjTF"
$ws.Range("G3").Value = "E44dyaN4k3F3MPkH"

# Update the active selection, as left by the editor
$ws.Range("N24").Select()
